# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to reflect the latest scraped counts.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 92
$wsExhibit.Range("F3").Value = 391
$wsExhibit.Range("F4").Value = 4976
$wsExhibit.Range("F5").Value = 35
$wsExhibit.Range("F6").Value = 34

# --- Sheet "全部类型" ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 92
$wsAll.Range("F3").Value = 391
$wsAll.Range("F4").Value = 4976
$wsAll.Range("F6").Value = 35
$wsAll.Range("F7").Value = 34
